$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 28 data (new site "Piney 24", duplicate lat/long of row 17 / "Piney 7" EPCHC dup)
$ws.Range("A28").Value = 27.626
$ws.Range("B28").Value = -82.5915
$ws.Range("C28").Value = "Piney 24"

# Apply border + right-aligned wrap-text style to A28:B28 (matches new cellXf s="3")
$borderRange = $ws.Range("A28:B28")
$borderRange.HorizontalAlignment = -4152  # xlRight
$borderRange.WrapText = $true
$borderRange.Borders.LineStyle = 1        # xlContinuous
$borderRange.Borders.Weight = -4138       # xlMedium
$borderRange.Borders.Color = 13421772

# Row 27 and 28 get a thick bottom border / slightly taller rows (ht=15, thickBot)
$ws.Rows.Item(27).RowHeight = 15
$ws.Rows.Item(28).RowHeight = 15
$ws.Range("A27:E28").Borders.Item(9).LineStyle = 1   # xlEdgeBottom? use explicit below
$ws.Range("A27:E27").Borders.Item(9).Weight = -4138

$ws.PageSetup.Orientation = 1  # xlPortrait

# Update view
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E23").Select()
